$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Urls")
$ws.Activate()

$ws.Range("A3").Value = "Employee List Page"
$ws.Range("B3").Value = "https://demo.actitime.com/administration/userlist.do"
$ws.Range("A4").Value = "Time Track"
$ws.Range("B4").Value = "https://demo.actitime.com/user/view_tt.do"

$ws.Columns.Item(1).ColumnWidth = 28.3
$ws.Columns.Item(2).ColumnWidth = 53.96
$ws.Columns.Item(5).ColumnWidth = 24.64

$ws.Range("B7").Select()
